$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rearranges the data (Fecha, Variedad, Calidad, Volumen, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg/unidad) across rows 2,3,4,5,6,7,8,10,12, while rows 1 (header), 9, 11
# and 13 stay untouched. Capture the "before" values first, then write the
# "after" values, so that reads always see original data regardless of
# write order.

$rows = @(2, 3, 4, 5, 6, 7, 8, 10, 12)

$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# destination row -> source row (the row whose original contents now land here)
$mapping = @{
    2  = 4
    3  = 7
    4  = 5
    5  = 8
    6  = 12
    7  = 3
    8  = 2
    10 = 6
    12 = 10
}

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $data = $before[$src]

    $ws.Cells.Item($dest, 4).Value2 = $data.D
    $ws.Cells.Item($dest, 11).Value2 = $data.K
    $ws.Cells.Item($dest, 12).Value2 = $data.L
    $ws.Cells.Item($dest, 13).Value2 = $data.M
    $ws.Cells.Item($dest, 14).Value2 = $data.N
    $ws.Cells.Item($dest, 15).Value2 = $data.O
    $ws.Cells.Item($dest, 16).Value2 = $data.P
    $ws.Cells.Item($dest, 17).Value2 = $data.Q
    $ws.Cells.Item($dest, 18).Value2 = $data.R
    $ws.Cells.Item($dest, 19).Value2 = $data.S
    $ws.Cells.Item($dest, 20).Value2 = $data.T
}
